$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -1
    4  = 1
    5  = -1
    6  = 2
    7  = 1
    8  = -1
    10 = 1
    11 = -7
    12 = 3
    13 = -4
    14 = -1
    16 = -2
    17 = 6
    18 = -3
    19 = 4
    20 = -4
    21 = 1
    22 = -7
    23 = 1
    24 = -4
    25 = 3
    26 = -1
    28 = 1
    29 = 1
    31 = 6
    32 = 2
    33 = -1
    34 = -3
    35 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
